$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7, pushing existing rows 7-21 down to 8-22
$ws.Rows.Item(7).Insert()

# Populate the new row 7 with the new weekly price record
$ws.Range("A7").Value = 10
$ws.Range("B7").Value = "Vega Modelo de Temuco"
$ws.Range("C7").Value = "La Araucanía"
$ws.Range("D7").Value = 45222
$ws.Range("E7").Value = 9
$ws.Range("F7").Value = "Fruta"
$ws.Range("G7").Value = 100104
$ws.Range("H7").Value = "Frutos de pepita"
$ws.Range("I7").Value = 100104004
$ws.Range("J7").Value = "Níspero"
$ws.Range("K7").Value = "Californiana(o)"
$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 25
$ws.Range("N7").Value = 28000
$ws.Range("O7").Value = 28000
$ws.Range("P7").Value = 28000
$ws.Range("Q7").Value = "$/bandeja 10 kilos"
$ws.Range("R7").Value = "Provincia de Quillota"
$ws.Range("S7").Value = 2800
$ws.Range("T7").Value = 10
